$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = "52.118.48"
$ws.Range("E2").Value = "  +0.85%  "
$rng.ClearFormats()

$rng = $ws.Range("D3:E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = "2.873.34"
$ws.Range("E3").Value = "  +3.12%  "
$rng.ClearFormats()

$rng = $ws.Range("E4")
$rng.NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$rng.ClearFormats()

$rng = $ws.Range("D5:E5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = "349.61"
$ws.Range("E5").Value = "  -0.91%  "
$rng.ClearFormats()

$rng = $ws.Range("D6:E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = "112.48"
$ws.Range("E6").Value = "  +3.24%  "
$rng.ClearFormats()

$rng = $ws.Range("D7:E7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = "0.555"
$ws.Range("E7").Value = "  +0.95%  "
$rng.ClearFormats()

$rng = $ws.Range("D8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$rng.ClearFormats()

$rng = $ws.Range("E9")
$rng.NumberFormat = "@"
$ws.Range("E9").Value = "  +1.80%  "
$rng.ClearFormats()

$rng = $ws.Range("D10:E10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = "40.23"
$ws.Range("E10").Value = "  +1.22%  "
$rng.ClearFormats()

$rng = $ws.Range("E11")
$rng.NumberFormat = "@"
$ws.Range("E11").Value = "  -0.44%  "
$rng.ClearFormats()

$rng = $ws.Range("E12")
$rng.NumberFormat = "@"
$ws.Range("E12").Value = "  +1.58%  "
$rng.ClearFormats()

$rng = $ws.Range("D13:E13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = "20.10"
$ws.Range("E13").Value = "  +0.32%  "
$rng.ClearFormats()

$rng = $ws.Range("D14:E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = "7.86"
$ws.Range("E14").Value = "  +2.20%  "
$rng.ClearFormats()

$rng = $ws.Range("D15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = "3.325.84"
$rng.ClearFormats()

$rng = $ws.Range("D16:E16")
$rng.NumberFormat = "@"
$ws.Range("D16").Value = "0.999"
$ws.Range("E16").Value = "  +7.54%  "
$rng.ClearFormats()

$rng = $ws.Range("D17:E17")
$rng.NumberFormat = "@"
$ws.Range("D17").Value = "2.861.83"
$ws.Range("E17").Value = "  +1.72%  "
$rng.ClearFormats()

$rng = $ws.Range("D18:E18")
$rng.NumberFormat = "@"
$ws.Range("D18").Value = "52.129.90"
$ws.Range("E18").Value = "  +0.84%  "
$rng.ClearFormats()

$rng = $ws.Range("E19")
$rng.NumberFormat = "@"
$ws.Range("E19").Value = "  +5.63%  "
$rng.ClearFormats()

$rng = $ws.Range("E20")
$rng.NumberFormat = "@"
$ws.Range("E20").Value = "  -1.52%  "
$rng.ClearFormats()

$rng = $ws.Range("D21:E21")
$rng.NumberFormat = "@"
$ws.Range("D21").Value = "13.63"
$ws.Range("E21").Value = "  +3.42%  "
$rng.ClearFormats()

$rng = $ws.Range("E22")
$rng.NumberFormat = "@"
$ws.Range("E22").Value = "  +1.11%  "
$rng.ClearFormats()

$rng = $ws.Range("D23:E23")
$rng.NumberFormat = "@"
$ws.Range("D23").Value = "70.69"
$ws.Range("E23").Value = "  +1.03%  "
$rng.ClearFormats()

$rng = $ws.Range("D24:E24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = "270.23"
$ws.Range("E24").Value = "  +1.03%  "
$rng.ClearFormats()

$rng = $ws.Range("D25:E25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = "2.78"
$ws.Range("E25").Value = "  +1.60%  "
$rng.ClearFormats()

$rng = $ws.Range("D26:E26")
$rng.NumberFormat = "@"
$ws.Range("D26").Value = "26.52"
$ws.Range("E26").Value = "  +1.81%  "
$rng.ClearFormats()

$rng = $ws.Range("E27")
$rng.NumberFormat = "@"
$ws.Range("E27").Value = "  +0.09%  "
$rng.ClearFormats()

$rng = $ws.Range("E28")
$rng.NumberFormat = "@"
$ws.Range("E28").Value = "  -0.51%  "
$rng.ClearFormats()

$rng = $ws.Range("D29:E29")
$rng.NumberFormat = "@"
$ws.Range("D29").Value = "10.58"
$ws.Range("E29").Value = "  +2.95%  "
$rng.ClearFormats()

$rng = $ws.Range("E30")
$rng.NumberFormat = "@"
$ws.Range("E30").Value = "  +2.72%  "
$rng.ClearFormats()

$rng = $ws.Range("D31:E31")
$rng.NumberFormat = "@"
$ws.Range("D31").Value = "6.27"
$ws.Range("E31").Value = "  +1.24%  "
$rng.ClearFormats()

$rng = $ws.Range("E32")
$rng.NumberFormat = "@"
$ws.Range("E32").Value = "  +1.14%  "
$rng.ClearFormats()

$rng = $ws.Range("D33:E33")
$rng.NumberFormat = "@"
$ws.Range("D33").Value = "5.82"
$ws.Range("E33").Value = "  +2.50%  "
$rng.ClearFormats()

$rng = $ws.Range("D34:E34")
$rng.NumberFormat = "@"
$ws.Range("D34").Value = "0.0455"
$ws.Range("E34").Value = "  -0.03%  "
$rng.ClearFormats()

$rng = $ws.Range("E35")
$rng.NumberFormat = "@"
$ws.Range("E35").Value = "  +7.57%  "
$rng.ClearFormats()

$rng = $ws.Range("E36")
$rng.NumberFormat = "@"
$ws.Range("E36").Value = "  -0.06%  "
$rng.ClearFormats()

$rng = $ws.Range("E37")
$rng.NumberFormat = "@"
$ws.Range("E37").Value = "  -15.73%  "
$rng.ClearFormats()

$rng = $ws.Range("D38:E38")
$rng.NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("E38").Value = "  +6.59%  "
$rng.ClearFormats()

$rng = $ws.Range("D39:E39")
$rng.NumberFormat = "@"
$ws.Range("D39").Value = "18.65"
$ws.Range("E39").Value = "  +0.90%  "
$rng.ClearFormats()

$rng = $ws.Range("E40")
$rng.NumberFormat = "@"
$ws.Range("E40").Value = "  +3.35%  "
$rng.ClearFormats()

$rng = $ws.Range("E41")
$rng.NumberFormat = "@"
$ws.Range("E41").Value = "  +4.12%  "
$rng.ClearFormats()

$rng = $ws.Range("E42")
$rng.NumberFormat = "@"
$ws.Range("E42").Value = "  +1.52%  "
$rng.ClearFormats()

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$rng = $ws.Range("D43:E43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = "22.66"
$ws.Range("E43").Value = "  +2.19%  "
$rng.ClearFormats()

$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$rng = $ws.Range("D44:E44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = "122.12"
$ws.Range("E44").Value = "  +1.26%  "
$rng.ClearFormats()

$rng = $ws.Range("D45:E45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = "2.21"
$ws.Range("E45").Value = "  +1.55%  "
$rng.ClearFormats()

$rng = $ws.Range("E46")
$rng.NumberFormat = "@"
$ws.Range("E46").Value = "  +5.42%  "
$rng.ClearFormats()

$rng = $ws.Range("D47:E47")
$rng.NumberFormat = "@"
$ws.Range("D47").Value = "2.173.05"
$ws.Range("E47").Value = "  +2.05%  "
$rng.ClearFormats()

$rng = $ws.Range("E48")
$rng.NumberFormat = "@"
$ws.Range("E48").Value = "  +6.36%  "
$rng.ClearFormats()

$rng = $ws.Range("E49")
$rng.NumberFormat = "@"
$ws.Range("E49").Value = "  +10.71%  "
$rng.ClearFormats()

$rng = $ws.Range("D50:E50")
$rng.NumberFormat = "@"
$ws.Range("D50").Value = "0.963"
$ws.Range("E50").Value = "  +6.20%  "
$rng.ClearFormats()

$rng = $ws.Range("E51")
$rng.NumberFormat = "@"
$ws.Range("E51").Value = "  +12.87%  "
$rng.ClearFormats()
